# Update the "想去人数" (F column) figures for the exhibitions sheets that
# changed between crawl runs ("展览" and "全部类型" contain the same rows).

$wb = $excel.ActiveWorkbook

# Row number (in both sheets) -> new value for column F
$updates = @{
    5  = 818
    6  = 271
    7  = 6539
    9  = 72
    10 = 110
    11 = 77
    14 = 15
    15 = 208
    16 = 528
    17 = 51
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
